$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) values per diff.
# Some Price values are plain decimal numbers (e.g. "582.97"); force
# those specific cells to Text format first so Excel keeps the exact
# literal string (preserving trailing zeros, etc.) instead of
# re-interpreting them as a number.

$ws.Range('D2').Value = '66.062.05'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '3.330.28'
$ws.Range('E3').Value = '  -0.33%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '582.97'
$ws.Range('E5').Value = '  +3.36%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '184.62'
$ws.Range('E6').Value = '  -2.54%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '3.326.83'
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.575'
$ws.Range('E9').Value = '  -2.63%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.180'
$ws.Range('E10').Value = '  -2.99%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.580'
$ws.Range('E11').Value = '  -2.06%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '46.93'
$ws.Range('E12').Value = '  -2.20%  '
$ws.Range('E13').Value = '  -2.04%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '664.94'
$ws.Range('E14').Value = '  +9.51%  '
$ws.Range('D15').Value = '3.868.52'
$ws.Range('E15').Value = '  -0.14%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '8.49'
$ws.Range('E16').Value = '  -2.66%  '
$ws.Range('D17').Value = '66.323.54'
$ws.Range('E17').Value = '  -0.20%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.117'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '17.84'
$ws.Range('E19').Value = '  -1.75%  '
$ws.Range('D20').Value = '3.333.54'
$ws.Range('E20').Value = '  -0.48%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.06'
$ws.Range('E21').Value = '  -1.27%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.896'
$ws.Range('E22').Value = '  -2.43%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '17.67'
$ws.Range('E23').Value = '  -6.17%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '102.07'
$ws.Range('E24').Value = '  +1.10%  '
$ws.Range('E25').Value = '  -3.03%  '
$ws.Range('E26').Value = '  -1.38%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.77'
$ws.Range('E27').Value = '  -0.39%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.42'
$ws.Range('E28').Value = '  -3.46%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '32.16'
$ws.Range('E29').Value = '  +4.97%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.47'
$ws.Range('E30').Value = '  -2.80%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.85'
$ws.Range('E31').Value = '  -0.20%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '610.48'
$ws.Range('E32').Value = '  +5.91%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.86'
$ws.Range('E33').Value = '  -1.85%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '11.08'
$ws.Range('E34').Value = '  -0.92%  '
$ws.Range('D35').Value = '3.857.40'
$ws.Range('E35').Value = '  +3.91%  '
$ws.Range('E36').Value = '  -1.26%  '
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '56.06'
$ws.Range('E38').Value = '  -2.27%  '
$ws.Range('E39').Value = '  -2.39%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.67'
$ws.Range('E40').Value = '  -1.86%  '
$ws.Range('D41').Value = '0.0₃0698'
$ws.Range('E41').Value = '  -4.89%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.41'
$ws.Range('E44').Value = '  +2.36%  '
$ws.Range('E45').Value = '  -3.14%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0415'
$ws.Range('E46').Value = '  -3.03%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.97'
$ws.Range('E47').Value = '  -13.85%  '
$ws.Range('E48').Value = '  -2.02%  '
$ws.Range('E49').Value = '  +0.47%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.54'
$ws.Range('E50').Value = '  -3.02%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.30'
$ws.Range('E51').Value = '  +1.09%  '

# Rows 42 and 43: the Coin name / Link / Price / Volume content moved
# (InjectiveProtocol now ranks above Stacks) -- update each cell to its
# new content accordingly.
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '32.72'
$ws.Range('E42').Value = '  -4.45%  '

$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.18'
$ws.Range('E43').Value = '  -4.18%  '
